$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:132 down to 104:133
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new record
$ws.Cells.Item(103, 1).Value = 5
$ws.Cells.Item(103, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(103, 3).Value = "Maule"
$ws.Cells.Item(103, 4).Value = 45218
$ws.Cells.Item(103, 5).Value = 7
$ws.Cells.Item(103, 6).Value = 100112026
$ws.Cells.Item(103, 7).Value = "Haba"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 400
$ws.Cells.Item(103, 11).Value = 10000
$ws.Cells.Item(103, 12).Value = 10000
$ws.Cells.Item(103, 13).Value = 10000
$ws.Cells.Item(103, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(103, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(103, 16).Value = 400
$ws.Cells.Item(103, 17).Value = 25
$ws.Cells.Item(103, 18).Value = "Hortaliza"
